$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSE")

# Duplicate the I3:K57 block (superdataset-11) into M3:O57 first (before J/K
# get their real values) so the new M/N/O block keeps the same "still empty"
# placeholder state that I/J/K currently have.
$src = $ws.Range("I3:K57")
$dst = $ws.Range("M3")
$src.Copy($dst)

# Update the new (3rd) block's title to the superdataset-12 dataset name.
$ws.Range("M3").Value = "Random Forest-100 (superdataset-12.csv)"

# Now fill in the previously-empty J/K (train/test MSE) values for the
# superdataset-11 block with the computed MSE results.
$jValues = @(
  0.000070874200231351315,
  0.0000608228496746671,
  0.000074588663772628575,
  0.00006760046075251362,
  0.000068633019270915841,
  0.000066988159880599052,
  0.000064653252323796127,
  0.000066273225706452559,
  0.000067172450574515673,
  0.000067089918896920056,
  0.000065712575381257068,
  0.000054382150437259108,
  0.000055847883628644771,
  0.000065348132055882788,
  0.00006702346638093486,
  0.000065555251651702452,
  0.000074563313736757532,
  0.000070051379744719879,
  0.000068756943354046399,
  0.00007474099130033646,
  0.000071017033131468205,
  0.000063591248191232331,
  0.000062545507369813662,
  0.000058781862659222063,
  0.000063310616339789727,
  0.000060818280722478357,
  0.000066338056252489077,
  0.000066262802899896874,
  0.000071605677193336065,
  0.000063623174724913829,
  0.000063644385237117667,
  0.000060686476680340451,
  0.000053914538224675222,
  0.000069706177037798836,
  0.00005847183910739454,
  0.000056466837989239277,
  0.000068446815664639565,
  0.000062258534378979718,
  0.000067552309498778039,
  0.000066381842477899225,
  0.000060905529793870131,
  0.000068621935687083105,
  0.00006277847551087128,
  0.000075820951200089224,
  0.000072460149064246712,
  0.000064655488684721592,
  0.00006356089105021878,
  0.000054675835871279189,
  0.000063701963003330052,
  0.000062834655894844461
)

$kValues = @(
  0.0003490213426590818,
  0.00035510666884648883,
  0.00040006883167582788,
  0.00029238078911649759,
  0.00039455950933599818,
  0.00026907721090905107,
  0.00096792392962882865,
  0.00028248852689565093,
  0.00035445521812474249,
  0.00044786002296544531,
  0.00040671076009897981,
  0.00065965253084427386,
  0.00066381434295174609,
  0.00038840750102423779,
  0.00041043540166760571,
  0.00034393959287596862,
  0.00038020975746799328,
  0.00034675786902055121,
  0.00043034771290308409,
  0.00032695894120276402,
  0.00043996624548749342,
  0.000384834694093958,
  0.00050684055392800616,
  0.00067346392218381522,
  0.0004883072888603317,
  0.00061530994821935576,
  0.0002668765886594242,
  0.00035716331436262401,
  0.00032746710305392663,
  0.00057007224774734762,
  0.00040948914896583399,
  0.00060989355957949333,
  0.00072623218603174073,
  0.00037374602692075772,
  0.00050484134207395988,
  0.00052298323997494278,
  0.00030002404100116641,
  0.0004044777870093374,
  0.00036227537414660119,
  0.00035115267775240658,
  0.00083962418889170075,
  0.00027880557414886541,
  0.00036153405812675092,
  0.00041705979114878348,
  0.00035708036486013919,
  0.00030372382159606347,
  0.0004259217153795682,
  0.00053411667181080149,
  0.00044711552840592561,
  0.00050111826616056369
)

for ($i = 0; $i -lt 50; $i++) {
  $row = 5 + $i
  $ws.Cells.Item($row, 10).Value = $jValues[$i]
  $ws.Cells.Item($row, 11).Value = $kValues[$i]
}

# Column widths for the new N/O columns (closest values reachable through the
# pixel-quantized ColumnWidth property).
$ws.Range("N1").ColumnWidth = 11.91650390625
$ws.Range("O1").ColumnWidth = 13.91650390625

# Restore the active cell selection.
$ws.Range("Q6").Select()
